$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 4; this pushes the existing rows 4-9 down to
# rows 5-10, keeping their content/heights intact.
$ws.Rows.Item(4).Insert()

# Copy the formatting (per-cell styles/borders) of the row above so the new
# row matches the look of the rest of the table.
$ws.Range("A3:F3").Copy($ws.Range("A4:F4"))

# Populate the newly inserted row with the "حذف مشتری" (Delete customer) entry.
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "حذف مشتری"
$ws.Range("C4").Value = 20
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "لاگین-جستجو بر اساس کد مشتری-حذف مشتری جستجو شده-کلیک روی دکمه ی ثبت"
$ws.Range("F4").ClearContents()
$ws.Rows.Item(4).RowHeight = 61.2

# The "شماره" numbering column is a plain 1..9 sequence; re-stamp it for every
# data row so it stays sequential after the insertion.
For ($i = 0; $i -le 8; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i + 1
}

# Match the selection left behind by the author (cell E4 active).
$ws.Range("E4").Select()
